$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$hl = $ws.Range("C4").Hyperlinks.Item(1)
$hl.TextToDisplay = ""
